$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: " week2" -> " week 2" in the title paragraph, and move the
# "_GoBack" bookmark (normally left behind by Word at the last edit
# location) to sit right after that run, inside the same paragraph.
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$text1 = $p1.Range.Text
$relStart = $text1.IndexOf(" week2")
if ($relStart -lt 0) {
    throw "Could not locate ' week2' in the first paragraph."
}
$absStart = $p1.Range.Start + $relStart
$absEnd = $absStart + (" week2").Length

$targetRange = $d.Range($absStart, $absEnd)

$weekRunXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
              '<w:r w:rsidR="00F47CCC">' +
                '<w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>' +
                '<w:t xml:space="preserve"> week 2</w:t>' +
              '</w:r>' +
              '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
              '<w:bookmarkEnd w:id="0"/>' +
              '</w:p>'

$targetRange.InsertXML($weekRunXml)

# ---------------------------------------------------------------------
# Edit 2: remove the old "_GoBack" bookmark that used to sit at the
# start of the "*Control statement -" paragraph.
# ---------------------------------------------------------------------
$controlPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.StartsWith("*Control statement")) {
        $controlPara = $candidate
        break
    }
}
if ($null -eq $controlPara) {
    throw "Could not locate the '*Control statement' paragraph."
}

$controlParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' +
                   'w:rsidR="008C68DD" w:rsidRDefault="008C68DD" w:rsidP="008C68DD">' +
                     '<w:r><w:t>*Control statement &#8211;</w:t></w:r>' +
                   '</w:p>'

$controlPara.Range.InsertXML($controlParaXml)
